# Add SVR parameter loading support: new columns K, L, M on Sheet1
# holding svr_kernel_scale, svr_epsilon and svr_box_constraint header/values,
# matching the "pred_par" structure used elsewhere in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New SVR parameter columns -------------------------------------------------
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.2
$ws.Range("M2").Value = 20

# --- Drop redundant explicit "Normal" styling on cells that only had the
# duplicate default cell style applied (cosmetic cleanup that happened when
# the workbook was resaved). ----------------------------------------------------
$normalCells = @(
    "B1","C1","H1","H2",
    "A5","A6","A8","A11",
    "D14","E14","F14","G14","I14",
    "D15","E15","F15","G15","I15",
    "D16","E16","F16","G16","I16"
)
foreach ($cellRef in $normalCells) {
    $ws.Range($cellRef).Style = "Normal"
}

# Row 7 had an explicit custom row format; clear it back to the sheet default.
$ws.Rows(7).ClearFormats()

# --- Selection moved to J9 in the resaved workbook -----------------------------
$ws.Range("J9").Select() | Out-Null
